$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextCell "D2" '63.120.15'
Set-TextCell "E2" '  -1.11%  '
Set-TextCell "D3" '3.081.68'
Set-TextCell "E3" '  +0.55%  '
Set-TextCell "E4" '  +0.05%  '
Set-TextCell "D5" '553.36'
Set-TextCell "E5" '  +0.79%  '
Set-TextCell "D6" '136.39'
Set-TextCell "E6" '  -3.46%  '
Set-TextCell "E7" '  -0.04%  '
Set-TextCell "D8" '3.073.73'
Set-TextCell "E8" '  +0.47%  '
Set-TextCell "D9" '0.493'
Set-TextCell "E9" '  +1.44%  '
Set-TextCell "E10" '  +2.49%  '
Set-TextCell "D11" '0.158'
Set-TextCell "E11" '  +5.59%  '
Set-TextCell "D12" '0.452'
Set-TextCell "E12" '  +1.56%  '
Set-TextCell "D13" '34.82'
Set-TextCell "E13" '  -1.14%  '
Set-TextCell "D14" '0.0000216'
Set-TextCell "E14" '  +1.18%  '
Set-TextCell "D15" '3.575.17'
Set-TextCell "E15" '  +0.61%  '
Set-TextCell "D16" '63.177.86'
Set-TextCell "E16" '  -1.13%  '
Set-TextCell "E17" '  +0.03%  '
Set-TextCell "D18" '3.084.16'
Set-TextCell "E18" '  +0.56%  '
Set-TextCell "D19" '499.82'
Set-TextCell "E19" '  +3.23%  '
Set-TextCell "D20" '6.63'
Set-TextCell "E20" '  +1.72%  '
Set-TextCell "D21" '13.45'
Set-TextCell "E21" '  -0.13%  '
Set-TextCell "D22" '0.702'
Set-TextCell "E22" '  +4.07%  '
Set-TextCell "D23" '7.24'
Set-TextCell "E23" '  +1.29%  '
Set-TextCell "D24" '77.76'
Set-TextCell "E24" '  +1.27%  '
Set-TextCell "D25" '12.26'
Set-TextCell "E25" '  -0.15%  '
Set-TextCell "E26" '  +0.13%  '
Set-TextCell "D27" '2.75'
Set-TextCell "E27" '  +2.40%  '
Set-TextCell "D28" '8.12'
Set-TextCell "E28" '  -0.40%  '
Set-TextCell "D29" '2.00'
Set-TextCell "E29" '  -2.43%  '
Set-TextCell "D30" '1.00'
Set-TextCell "E30" '  -0.01%  '
Set-TextCell "D31" '26.16'
Set-TextCell "E31" '  +2.57%  '
Set-TextCell "D32" '2.50'
Set-TextCell "E32" '  -4.12%  '
Set-TextCell "E33" '  -1.46%  '
Set-TextCell "D34" '59.16'
Set-TextCell "E34" '  +13.86%  '
Set-TextCell "D35" '529.00'
Set-TextCell "E35" '  -8.22%  '
Set-TextCell "D36" '5.86'
Set-TextCell "E36" '  +1.16%  '
Set-TextCell "D37" '5.13'
Set-TextCell "E37" '  -2.54%  '
Set-TextCell "D38" '0.0409'
Set-TextCell "E38" '  +3.28%  '
Set-TextCell "D39" '3.047.40'
Set-TextCell "E39" '  +2.40%  '
Set-TextCell "D40" '0.0788'
Set-TextCell "E40" '  +1.05%  '
Set-TextCell "E41" '  +3.09%  '
Set-TextCell "D42" '8.04'
Set-TextCell "E42" '  -0.44%  '
Set-TextCell "D43" '2.61'
Set-TextCell "E43" '  -6.64%  '
Set-TextCell "D44" '0.252'
Set-TextCell "E44" '  +4.40%  '
Set-TextCell "E45" '  +0.01%  '
Set-TextCell "D46" '2.05'
Set-TextCell "E46" '  -0.65%  '
Set-TextCell "D47" '120.04'
Set-TextCell "E47" '  +2.15%  '
Set-TextCell "E48" '  -0.16%  '
Set-TextCell "D49" '23.62'
Set-TextCell "E49" '  -4.43%  '
Set-TextCell "D50" '0.0₃0492'
Set-TextCell "E50" '  -4.58%  '
Set-TextCell "D51" '2.33'
Set-TextCell "E51" '  +63.46%  '
